$wb = $excel.ActiveWorkbook

# --- 1. Remove the obsolete "Sheet" row (row 16) from the
#        optimization_parameters sheet. This row held the header label
#        "Sheet" together with values 3 / 4 that are no longer needed;
#        deleting it shifts the "simulation_timepoints" row up from 17 to 16
#        and drops the now-unused "Sheet" shared string.
$wsParams = $wb.Worksheets.Item("optimization_parameters")
$wsParams.Rows.Item(16).Delete() | Out-Null

# Leave the selection on the row that now occupies position 16 (the row
# that used to be row 17), matching how Excel leaves the selection after
# deleting a whole row.
$wsParams.Rows.Item(16).Select() | Out-Null

# --- 2. Switch the active sheet to optimization_diagnostics (the last
#        sheet in the workbook), as the last step of the audit.
$wsDiag = $wb.Worksheets.Item("optimization_diagnostics")
$wsDiag.Activate() | Out-Null
